# Auto-applies the numeric corrections described in the commit diff
# (re-synced price/profit figures for several Leve rows across sheets).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 4013.377
$ws.Range("I138").Value = 1668.75
$ws.Range("J138").Value = 7389.64
$ws.Range("K138").Value = 5006.25
$ws.Range("L138").Value = 22168.92
$ws.Range("M138").Value = 133.75
$ws.Range("N138").Value = -32448.92

$ws.Range("H141").Value = 786933.5
$ws.Range("I141").Value = 2493.125
$ws.Range("J141").Value = 1484213.8
$ws.Range("K141").Value = 7479.375
$ws.Range("L141").Value = 4452641.4
$ws.Range("M141").Value = -2299.375
$ws.Range("N141").Value = -4463001.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 16668601
$ws.Range("I2").Value = 125000210
$ws.Range("J2").Value = 2200
$ws.Range("K2").Value = 125000210
$ws.Range("L2").Value = 2200
$ws.Range("M2").Value = -125000097
$ws.Range("N2").Value = -2426

$ws.Range("H32").Value = 3687.592
$ws.Range("I32").Value = 2828.875
$ws.Range("J32").Value = 8267.416999999999
$ws.Range("K32").Value = 2828.875
$ws.Range("L32").Value = 8267.416999999999
$ws.Range("M32").Value = -2541.875
$ws.Range("N32").Value = -8841.416999999999

$ws.Range("H45").Value = 1712.1471
$ws.Range("I45").Value = 1074.9642
$ws.Range("K45").Value = 1074.9642
$ws.Range("M45").Value = -697.9641999999999

$ws.Range("H61").Value = 4349.5
$ws.Range("I61").Value = 1741.7142
$ws.Range("J61").Value = 6377.778
$ws.Range("K61").Value = 1741.7142
$ws.Range("L61").Value = 6377.778
$ws.Range("M61").Value = -1529.7142
$ws.Range("N61").Value = -6801.778

$ws.Range("H88").Value = 3480
$ws.Range("I88").Value = 1975
$ws.Range("J88").Value = 9500
$ws.Range("K88").Value = 1975
$ws.Range("L88").Value = 9500
$ws.Range("M88").Value = -1569
$ws.Range("N88").Value = -10312

$ws.Range("H91").Value = 3480
$ws.Range("I91").Value = 1975
$ws.Range("J91").Value = 9500
$ws.Range("K91").Value = 1975
$ws.Range("L91").Value = 9500
$ws.Range("M91").Value = -571
$ws.Range("N91").Value = -12308

$ws.Range("H110").Value = 1298.8889
$ws.Range("I110").Value = 725
$ws.Range("J110").Value = 2446.6667
$ws.Range("K110").Value = 725
$ws.Range("L110").Value = 2446.6667
$ws.Range("M110").Value = 1320
$ws.Range("N110").Value = -6536.6667

$ws.Range("H116").Value = 16668601
$ws.Range("I116").Value = 125000210
$ws.Range("J116").Value = 2200
$ws.Range("K116").Value = 125000210
$ws.Range("L116").Value = 2200
$ws.Range("M116").Value = -124997916
$ws.Range("N116").Value = -6788

$ws.Range("H132").Value = 17546892
$ws.Range("I132").Value = 22225042
$ws.Range("K132").Value = 66675126
$ws.Range("M132").Value = -66672596

$ws.Range("H136").Value = 4349.5
$ws.Range("I136").Value = 1741.7142
$ws.Range("J136").Value = 6377.778
$ws.Range("K136").Value = 5225.142599999999
$ws.Range("L136").Value = 19133.334
$ws.Range("M136").Value = -2675.142599999999
$ws.Range("N136").Value = -24233.334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 16668601
$ws.Range("I3").Value = 125000210
$ws.Range("J3").Value = 2200
$ws.Range("K3").Value = 125000210
$ws.Range("L3").Value = 2200
$ws.Range("M3").Value = -125000096
$ws.Range("N3").Value = -2428

$ws.Range("H86").Value = 1438985.9
$ws.Range("I86").Value = 1917058.9
$ws.Range("J86").Value = 4766.6665
$ws.Range("K86").Value = 1917058.9
$ws.Range("L86").Value = 4766.6665
$ws.Range("M86").Value = -1915935.9
$ws.Range("N86").Value = -7012.6665

$ws.Range("H89").Value = 1438985.9
$ws.Range("I89").Value = 1917058.9
$ws.Range("J89").Value = 4766.6665
$ws.Range("K89").Value = 9585294.5
$ws.Range("L89").Value = 23833.3325
$ws.Range("M89").Value = -9579678.5
$ws.Range("N89").Value = -35065.3325

$ws.Range("H99").Value = 1623.4783
$ws.Range("I99").Value = 1044.9286
$ws.Range("K99").Value = 1044.9286
$ws.Range("M99").Value = 453.0714

$ws.Range("H107").Value = 4300
$ws.Range("I107").Value = 1400
$ws.Range("J107").Value = 5750
$ws.Range("K107").Value = 1400
$ws.Range("L107").Value = 5750
$ws.Range("M107").Value = 520
$ws.Range("N107").Value = -9590

$ws.Range("H134").Value = 3911.647
$ws.Range("I134").Value = 2653.7693
$ws.Range("K134").Value = 7961.3079
$ws.Range("M134").Value = -5426.3079

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 4274.5
$ws.Range("I16").Value = 2100
$ws.Range("J16").Value = 4999.3335
$ws.Range("K16").Value = 2100
$ws.Range("L16").Value = 4999.3335
$ws.Range("M16").Value = -1813
$ws.Range("N16").Value = -5573.3335

$ws.Range("H105").Value = 10650
$ws.Range("I105").Value = 10866.667
$ws.Range("J105").Value = 10000
$ws.Range("K105").Value = 10866.667
$ws.Range("L105").Value = 10000
$ws.Range("M105").Value = -9119.666999999999
$ws.Range("N105").Value = -13494

$ws.Range("H113").Value = 4274.5
$ws.Range("I113").Value = 2100
$ws.Range("J113").Value = 4999.3335
$ws.Range("K113").Value = 2100
$ws.Range("L113").Value = 4999.3335
$ws.Range("M113").Value = 70
$ws.Range("N113").Value = -9339.333500000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 3841.4285
$ws.Range("I70").Value = 3809.0908
$ws.Range("J70").Value = 3960
$ws.Range("K70").Value = 3809.0908
$ws.Range("L70").Value = 3960
$ws.Range("M70").Value = -3539.0908
$ws.Range("N70").Value = -4500

$ws.Range("H73").Value = 3841.4285
$ws.Range("I73").Value = 3809.0908
$ws.Range("J73").Value = 3960
$ws.Range("K73").Value = 3809.0908
$ws.Range("L73").Value = 3960
$ws.Range("M73").Value = -2873.0908
$ws.Range("N73").Value = -5832

$ws.Range("H132").Value = 43482140
$ws.Range("I132").Value = 83335430
$ws.Range("J132").Value = 5818
$ws.Range("K132").Value = 250006290
$ws.Range("L132").Value = 17454
$ws.Range("M132").Value = -250003760
$ws.Range("N132").Value = -22514

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("M64").ClearContents()

$ws.Range("H67").Value = 0
$ws.Range("I67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("M67").ClearContents()

$ws.Range("H95").Value = 30344
$ws.Range("J95").Value = 30344
$ws.Range("L95").Value = 30344
$ws.Range("N95").Value = -35836

$ws.Range("H100").Value = 2170.6924
$ws.Range("I100").Value = 1250
$ws.Range("J100").Value = 2959.8572
$ws.Range("K100").Value = 1250
$ws.Range("L100").Value = 2959.8572
$ws.Range("M100").Value = -709
$ws.Range("N100").Value = -4041.8572

$ws.Range("H136").Value = 1991.4667
$ws.Range("I136").Value = 1591.2222
$ws.Range("J136").Value = 2591.8333
$ws.Range("K136").Value = 4773.6666
$ws.Range("L136").Value = 7775.499899999999
$ws.Range("M136").Value = -2223.6666
$ws.Range("N136").Value = -12875.4999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 4556
$ws.Range("I29").Value = 4860
$ws.Range("J29").Value = 4100
$ws.Range("K29").Value = 4860
$ws.Range("L29").Value = 4100
$ws.Range("M29").Value = -4570
$ws.Range("N29").Value = -4680

$ws.Range("H97").Value = 32000
$ws.Range("J97").Value = 32000
$ws.Range("L97").Value = 32000
$ws.Range("N97").Value = -33982

$ws.Range("H122").Value = 2486.25
$ws.Range("I122").Value = 1952.3077
$ws.Range("K122").Value = 5856.9231
$ws.Range("M122").Value = -3406.9231

$ws.Range("H126").Value = 2430.9583
$ws.Range("I126").Value = 1964.4615
$ws.Range("J126").Value = 2982.2727
$ws.Range("K126").Value = 5893.3845
$ws.Range("L126").Value = 8946.8181
$ws.Range("M126").Value = -3423.3845
$ws.Range("N126").Value = -13886.8181

$ws.Range("H132").Value = 14816.122
$ws.Range("I132").Value = 2343.6365
$ws.Range("J132").Value = 29257.947
$ws.Range("K132").Value = 7030.9095
$ws.Range("L132").Value = 87773.841
$ws.Range("M132").Value = -4500.9095
$ws.Range("N132").Value = -92833.841

$ws.Range("H136").Value = 1974.4
$ws.Range("I136").Value = 872.5
$ws.Range("J136").Value = 2617.5
$ws.Range("K136").Value = 2617.5
$ws.Range("M136").Value = -67.5
